# RETC_062 and RETC_063 upload
# Adds a new "Sheet2" (placed after "Sheet1") with an Email / First Name /
# Last Name table, reproducing the rows previously held on Sheet1, and makes
# Sheet2 the active sheet/tab.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")

# New worksheet goes right after Sheet1, named "Sheet2" - this also makes
# it the active sheet (tabSelected/activeTab follow the newly added sheet).
$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "Sheet2"

# Write row 2 (manzoor) before the header row so the shared-string table
# picks up "manzoor" ahead of the header labels, matching source order.
$ws.Range("A2").Value = "manzoor"
$ws.Range("B2").Value = "manzoor mehadi"
$ws.Range("C2").Value = "manzoor mehadi"

$ws.Range("A1").Value = "Email"
$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Last Name"

$ws.Range("A3").Value = "alex@gmail.com"
$ws.Range("B3").Value = 76576
$ws.Range("C3").Value = "alex hales"

$ws.Range("A4").Value = "mariya@gmail.com"
$ws.Range("B4").Value = "mariya"
$ws.Range("C4").Value = 876876

# Mailto hyperlinks on column A (rows 2-4), mirroring Sheet1's pattern.
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:manzoor@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:alex@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:mariya@gmail.com")

# Match Sheet1's existing formatting (thin-bordered normal / hyperlink
# styles) instead of the style Hyperlinks.Add / plain writes picked up.
# Applied last so it wins over the auto-assigned hyperlink style.
$sheet1.Range("B1").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$ws.Range("B2:C4").PasteSpecial(-4122)

$sheet1.Range("A1").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)

$null = $ws.Range("A1").Select()
